$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.658.03"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "2.590.15"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.58%  "

$ws.Range("D9").Value = "2.594.66"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.14%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("D14").Value = "3.042.26"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").Value = "60.583.42"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.76%  "

$ws.Range("D18").Value = "2.593.50"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").Value = "0.0₃0846"
$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.14%  "

$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.856"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.19%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0557"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.84%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "

$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("D51").Value = "2.003.45"
$ws.Range("E51").Value = "  +0.06%  "
